{"js": "// Each table cell in this worksheet holds a '<three-digit>\u00d7<one-digit>=<product>'\n// multiplication fact as plain text. Replace every old fact with its new value,\n// matching on the exact original text (all 25 values in the sheet are unique).\nconst replacements = [\n  [\"513\u00d78=4104\", \"508\u00d74=2032\"],\n  [\"541\u00d73=1623\", \"398\u00d79=3582\"],\n  [\"308\u00d72=616\", \"967\u00d74=3868\"],\n  [\"765\u00d79=6885\", \"225\u00d72=450\"],\n  [\"539\u00d72=1078\", \"219\u00d72=438\"],\n  [\"686\u00d74=2744\", \"545\u00d73=1635\"],\n  [\"430\u00d74=1720\", \"794\u00d74=3176\"],\n  [\"696\u00d79=6264\", \"317\u00d76=1902\"],\n  [\"229\u00d74=916\", \"360\u00d75=1800\"],\n  [\"674\u00d77=4718\", \"993\u00d79=8937\"],\n  [\"397\u00d76=2382\", \"442\u00d79=3978\"],\n  [\"913\u00d78=7304\", \"257\u00d75=1285\"],\n  [\"453\u00d76=2718\", \"133\u00d72=266\"],\n  [\"358\u00d72=716\", \"322\u00d75=1610\"],\n  [\"426\u00d73=1278\", \"728\u00d76=4368\"],\n  [\"738\u00d77=5166\", \"846\u00d72=1692\"],\n  [\"712\u00d76=4272\", \"162\u00d73=486\"],\n  [\"813\u00d75=4065\", \"682\u00d73=2046\"],\n  [\"391\u00d75=1955\", \"920\u00d76=5520\"],\n  [\"466\u00d76=2796\", \"992\u00d74=3968\"],\n  [\"961\u00d73=2883\", \"621\u00d76=3726\"],\n  [\"929\u00d74=3716\", \"589\u00d76=3534\"],\n  [\"206\u00d78=1648\", \"134\u00d72=268\"],\n  [\"186\u00d78=1488\", \"962\u00d75=4810\"],\n  [\"766\u00d77=5362\", \"321\u00d77=2247\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit x one-digit multiplication answer in the\n# document's table cells with its new value, matched by exact old text.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"513\u00d78=4104\", \"508\u00d74=2032\"),\n  @(\"541\u00d73=1623\", \"398\u00d79=3582\"),\n  @(\"308\u00d72=616\", \"967\u00d74=3868\"),\n  @(\"765\u00d79=6885\", \"225\u00d72=450\"),\n  @(\"539\u00d72=1078\", \"219\u00d72=438\"),\n  @(\"686\u00d74=2744\", \"545\u00d73=1635\"),\n  @(\"430\u00d74=1720\", \"794\u00d74=3176\"),\n  @(\"696\u00d79=6264\", \"317\u00d76=1902\"),\n  @(\"229\u00d74=916\", \"360\u00d75=1800\"),\n  @(\"674\u00d77=4718\", \"993\u00d79=8937\"),\n  @(\"397\u00d76=2382\", \"442\u00d79=3978\"),\n  @(\"913\u00d78=7304\", \"257\u00d75=1285\"),\n  @(\"453\u00d76=2718\", \"133\u00d72=266\"),\n  @(\"358\u00d72=716\", \"322\u00d75=1610\"),\n  @(\"426\u00d73=1278\", \"728\u00d76=4368\"),\n  @(\"738\u00d77=5166\", \"846\u00d72=1692\"),\n  @(\"712\u00d76=4272\", \"162\u00d73=486\"),\n  @(\"813\u00d75=4065\", \"682\u00d73=2046\"),\n  @(\"391\u00d75=1955\", \"920\u00d76=5520\"),\n  @(\"466\u00d76=2796\", \"992\u00d74=3968\"),\n  @(\"961\u00d73=2883\", \"621\u00d76=3726\"),\n  @(\"929\u00d74=3716\", \"589\u00d76=3534\"),\n  @(\"206\u00d78=1648\", \"134\u00d72=268\"),\n  @(\"186\u00d78=1488\", \"962\u00d75=4810\"),\n  @(\"766\u00d77=5362\", \"321\u00d77=2247\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  # wdFindContinue = 1, wdReplaceAll = 2; MatchCase = $true for an exact match\n  $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"No match found for: $oldText\"\n  }\n}\n"}
